$wb = $excel.ActiveWorkbook

$wsBuy = $wb.Worksheets.Item("open_buy_orders")
$wsSell = $wb.Worksheets.Item("open_sell_orders")

# open_buy_orders: remove row 2 (ONIP37-KFKJF-LO62AZ / 1.9735)
$wsBuy.Rows.Item(2).Delete()

# open_sell_orders: row 2 becomes OAGYX4-QHUTU-OJT25D, remove row 3
$wsSell.Range("A2").Value = "OAGYX4-QHUTU-OJT25D"
$wsSell.Rows.Item(3).Delete()
